$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2706.8572
$ws.Range("I98").Value = 2601
$ws.Range("J98").Value = 3095
$ws.Range("K98").Value = 2601
$ws.Range("L98").Value = 3095
$ws.Range("M98").Value = -1103
$ws.Range("N98").Value = -6091
$ws.Range("H116").Value = 3150.2
$ws.Range("I116").Value = 2384.5386
$ws.Range("J116").Value = 4572.143
$ws.Range("K116").Value = 2384.5386
$ws.Range("L116").Value = 4572.143
$ws.Range("M116").Value = 1057.4614
$ws.Range("N116").Value = -11456.143
$ws.Range("H122").Value = 2706.8572
$ws.Range("I122").Value = 2601
$ws.Range("J122").Value = 3095
$ws.Range("K122").Value = 7803
$ws.Range("L122").Value = 9285
$ws.Range("M122").Value = -5353
$ws.Range("N122").Value = -14185
$ws.Range("H132").Value = 4674.8335
$ws.Range("I132").Value = 2343.4092
$ws.Range("J132").Value = 11086.25
$ws.Range("K132").Value = 7030.2276
$ws.Range("L132").Value = 33258.75
$ws.Range("M132").Value = -4500.2276
$ws.Range("N132").Value = -38318.75
$ws.Range("H135").Value = 115.888885
$ws.Range("I135").Value = 80.375
$ws.Range("K135").Value = 723.375
$ws.Range("M135").Value = 1811.625
$ws.Range("H137").Value = 1413.2
$ws.Range("I137").Value = 1280.1818
$ws.Range("J137").Value = 1575.7778
$ws.Range("K137").Value = 3840.5454
$ws.Range("L137").Value = 4727.3334
$ws.Range("M137").Value = -1290.5454
$ws.Range("N137").Value = -9827.3334
$ws.Range("H138").Value = 589485.4399999999
$ws.Range("J138").Value = 792706.7
$ws.Range("L138").Value = 2378120.1
$ws.Range("N138").Value = -2388400.1
$ws.Range("H140").Value = 48766.668
$ws.Range("J140").Value = 48766.668
$ws.Range("L140").Value = 48766.668
$ws.Range("N140").Value = -59126.668
# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3289.5715
$ws.Range("I32").Value = 3351.5122
$ws.Range("K32").Value = 3351.5122
$ws.Range("M32").Value = -3064.5122
$ws.Range("H45").Value = 1635.6111
$ws.Range("I45").Value = 1601.8
$ws.Range("J45").Value = 1804.6666
$ws.Range("K45").Value = 1601.8
$ws.Range("L45").Value = 1804.6666
$ws.Range("M45").Value = -1224.8
$ws.Range("N45").Value = -2558.6666
$ws.Range("H61").Value = 2299
$ws.Range("I61").Value = 1626.6666
$ws.Range("J61").Value = 2971.3333
$ws.Range("K61").Value = 1626.6666
$ws.Range("L61").Value = 2971.3333
$ws.Range("M61").Value = -1414.6666
$ws.Range("N61").Value = -3395.3333
$ws.Range("H74").Value = 781.069
$ws.Range("I74").Value = 769.2
$ws.Range("K74").Value = 769.2
$ws.Range("M74").Value = 104.8
$ws.Range("H77").Value = 781.069
$ws.Range("I77").Value = 769.2
$ws.Range("K77").Value = 3846
$ws.Range("M77").Value = 522
$ws.Range("H110").Value = 1166.4857
$ws.Range("I110").Value = 975.36664
$ws.Range("K110").Value = 975.36664
$ws.Range("M110").Value = 1069.63336
$ws.Range("H132").Value = 4178.9287
$ws.Range("I132").Value = 4250.7
$ws.Range("K132").Value = 12752.1
$ws.Range("M132").Value = -10222.1
$ws.Range("H136").Value = 2299
$ws.Range("I136").Value = 1626.6666
$ws.Range("J136").Value = 2971.3333
$ws.Range("K136").Value = 4879.9998
$ws.Range("L136").Value = 8913.999899999999
$ws.Range("M136").Value = -2329.9998
$ws.Range("N136").Value = -14013.9999
# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 25001204
$ws.Range("I99").Value = 27778944
$ws.Range("K99").Value = 27778944
$ws.Range("M99").Value = -27777446
$ws.Range("H105").Value = 72137560
$ws.Range("I105").Value = 100991176
$ws.Range("J105").Value = 3500
$ws.Range("K105").Value = 100991176
$ws.Range("L105").Value = 3500
$ws.Range("M105").Value = -100989429
$ws.Range("N105").Value = -6994
$ws.Range("H132").Value = 15000
$ws.Range("J132").Value = 15000
$ws.Range("L132").Value = 15000
$ws.Range("N132").Value = -25120
$ws.Range("H134").Value = 9729.538
$ws.Range("I134").Value = 2048.5
$ws.Range("K134").Value = 6145.5
$ws.Range("M134").Value = -3610.5
# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 888.9792
$ws.Range("J31").Value = 1440.909
$ws.Range("L31").Value = 1440.909
$ws.Range("N31").Value = -2030.909
$ws.Range("H34").Value = 888.9792
$ws.Range("J34").Value = 1440.909
$ws.Range("L34").Value = 1440.909
$ws.Range("N34").Value = -1844.909
$ws.Range("H58").Value = 744.4167
$ws.Range("I58").Value = 726.8570999999999
$ws.Range("J58").Value = 769
$ws.Range("K58").Value = 726.8570999999999
$ws.Range("L58").Value = 769
$ws.Range("M58").Value = -523.8570999999999
$ws.Range("N58").Value = -1175
$ws.Range("H132").Value = 5760.7036
$ws.Range("J132").Value = 3055.4546
$ws.Range("L132").Value = 9166.363799999999
$ws.Range("N132").Value = -14226.3638
$ws.Range("H135").Value = 66833.336
$ws.Range("J135").Value = 66833.336
$ws.Range("L135").Value = 66833.336
$ws.Range("N135").Value = -76973.336
$ws.Range("H136").Value = 744.4167
$ws.Range("I136").Value = 726.8570999999999
$ws.Range("J136").Value = 769
$ws.Range("K136").Value = 2180.5713
$ws.Range("L136").Value = 2307
$ws.Range("M136").Value = 369.4287000000004
$ws.Range("N136").Value = -7407
# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1712.52
$ws.Range("I5").Value = 1578.2778
$ws.Range("K5").Value = 4734.8334
$ws.Range("M5").Value = -4622.8334
$ws.Range("H12").Value = 796.6667
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 796.6667
$ws.Range("K12").Value = 0
$ws.Range("L12").ClearContents()
$ws.Range("M12").Value = 2390.0001
$ws.Range("N12").Value = -2736.0001
$ws.Range("H13").Value = 637.4
$ws.Range("I13").Value = 193.5
$ws.Range("J13").Value = 933.3333
$ws.Range("K13").Value = 580.5
$ws.Range("L13").Value = 2799.9999
$ws.Range("M13").Value = -412.5
$ws.Range("N13").Value = -3135.9999
$ws.Range("H98").Value = 750
$ws.Range("I98").Value = 126.666664
$ws.Range("J98").Value = 1996.6666
$ws.Range("K98").Value = 379.999992
$ws.Range("L98").Value = 5989.9998
$ws.Range("M98").Value = 1118.000008
$ws.Range("N98").Value = -8985.9998
$ws.Range("H127").Value = 1774.9
$ws.Range("J127").Value = 1774.9
$ws.Range("L127").Value = 5324.700000000001
$ws.Range("N127").Value = -15244.7
$ws.Range("H131").Value = 26317292
$ws.Range("I131").Value = 166667360
$ws.Range("J131").Value = 1654.3125
$ws.Range("K131").Value = 500002080
$ws.Range("L131").Value = 4962.9375
$ws.Range("M131").Value = -499997040
$ws.Range("N131").Value = -15042.9375
$ws.Range("H135").Value = 1712.52
$ws.Range("I135").Value = 1578.2778
$ws.Range("K135").Value = 14204.5002
$ws.Range("M135").Value = -11669.5002
# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3220.158
$ws.Range("I122").Value = 2599
$ws.Range("J122").Value = 3582.5
$ws.Range("K122").Value = 7797
$ws.Range("L122").Value = 10747.5
$ws.Range("M122").Value = -5347
$ws.Range("N122").Value = -15647.5
$ws.Range("H132").Value = 2690.0417
$ws.Range("I132").Value = 2304.6428
$ws.Range("K132").Value = 6913.928400000001
$ws.Range("M132").Value = -4383.928400000001
$ws.Range("H135").Value = 47749
$ws.Range("J135").Value = 47749
$ws.Range("L135").Value = 47749
$ws.Range("N135").Value = -57889
# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 429285.72
$ws.Range("I2").Value = 500000
$ws.Range("J2").Value = 302000
$ws.Range("K2").Value = 500000
$ws.Range("L2").Value = 302000
$ws.Range("M2").Value = -499888
$ws.Range("N2").Value = -302224
$ws.Range("H22").Value = 1414.2858
$ws.Range("I22").Value = 1345.5454
$ws.Range("K22").Value = 1345.5454
$ws.Range("M22").Value = -1050.5454
$ws.Range("H27").Value = 1414.2858
$ws.Range("I27").Value = 1345.5454
$ws.Range("K27").Value = 1345.5454
$ws.Range("M27").Value = -1238.5454
$ws.Range("H40").Value = 2355.375
$ws.Range("I40").Value = 2191.2856
$ws.Range("J40").Value = 3504
$ws.Range("K40").Value = 2191.2856
$ws.Range("L40").Value = 3504
$ws.Range("M40").Value = -2055.2856
$ws.Range("N40").Value = -3776
$ws.Range("H61").Value = 2441.4092
$ws.Range("I61").Value = 1941.9166
$ws.Range("J61").Value = 3040.8
$ws.Range("K61").Value = 1941.9166
$ws.Range("L61").Value = 3040.8
$ws.Range("M61").Value = -1739.9166
$ws.Range("N61").Value = -3444.8
$ws.Range("H113").Value = 2441.4092
$ws.Range("I113").Value = 1941.9166
$ws.Range("J113").Value = 3040.8
$ws.Range("K113").Value = 1941.9166
$ws.Range("L113").Value = 3040.8
$ws.Range("M113").Value = 228.0834
$ws.Range("N113").Value = -7380.8
# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 439.36
$ws.Range("I113").Value = 252.11765
$ws.Range("J113").Value = 837.25
$ws.Range("K113").Value = 756.35295
$ws.Range("L113").Value = 2511.75
$ws.Range("M113").Value = 1413.64705
$ws.Range("N113").Value = -6851.75
$ws.Range("H132").Value = 2664.45
$ws.Range("I132").Value = 2399.1428
$ws.Range("J132").Value = 3283.5
$ws.Range("K132").Value = 7197.428400000001
$ws.Range("L132").Value = 9850.5
$ws.Range("M132").Value = -4667.428400000001
$ws.Range("N132").Value = -14910.5
$ws.Range("H136").Value = 2159.6
$ws.Range("I136").Value = 1933
$ws.Range("K136").Value = 5799
$ws.Range("M136").Value = -3249
